$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.153579950332642
$ws.Range("B1").Value = 2.216450929641724
$ws.Range("C1").Value = 2.296318531036377
$ws.Range("D1").Value = 3.051946640014648
$ws.Range("E1").Value = 2.759232759475708
